$wb = $excel.ActiveWorkbook

# Overview sheet: row 3 corresponds to the 879e6e22-... file.
# Its status in both the zh-cn and de-de columns moves from
# "Ready for handoff" to "Handed back: in sync with en-US".
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 (879e6e22-... file) gets the same status update,
# plus a refreshed "Latest Handback DateTime" (column G).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("G3").Value = "2016-01-18 02:52:22"

# de-de sheet: row 3 (879e6e22-... file) gets the same status update,
# plus a refreshed "Latest Handback DateTime" (column G).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("G3").Value = "2016-01-18 02:52:39"
